$wb = $excel.ActiveWorkbook

# --- "Constant Samples" sheet: bump zoom to 201% before duplicating it,
#     so the new sheet inherits the same view setting. ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.Zoom = 201

# --- Duplicate it right after itself to become "Formula Samples". ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Formula Samples"

# Row 3 ("Decade"): constant value left as-is, but flagged as wrong since this
# sheet is meant to exercise formulas.
$ws2.Range("C3").Value = "> This cell will be considered as wrong because it inputs a constant value, not a formula."

# Row 7 ("Kilometer"): turn the constant into a formula, drop its old comment.
$ws2.Range("B7").Formula = "= B6 * 0.001"
$ws2.Range("C7").ClearContents()

# Row 8 ("Feet"): turn the constant into a (slightly off) formula, new comment.
$ws2.Range("B8").Formula = "= B6 * 3.28"
$ws2.Range("C8").Value = "> This will be considered as wrong due to incorrect formula."

# Row 4 ("Century"): turn the constant into a ROUNDUP formula, new comment.
$ws2.Range("B4").Formula = "= ROUNDUP(B2 * 0.01, 0)"
$ws2.Range("C4").Value = "> This cel will be considered as wrong since ROUNDUP is yet to be implemented."

# Row 9 ("Inch"): turn the constant into a formula; comment text unchanged.
$ws2.Range("B9").Formula = "=B6*39.37"
$ws2.Range("C9").Value = "> This will be considered as correct due to the alt_cells."

$ws2.Range("C5").Select()
